# fix format number harga input with js
# Rewrites the vehicle-asset header row to snake_case field keys and
# refreshes the sample data row (Avanza -> Xenia) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header labels -> snake_case field keys ---------------------
$ws.Range("A1").Value = "id_status_aset"
$ws.Range("B1").Value = "nama"
$ws.Range("C1").Value = "tanggal_inventarisir"
$ws.Range("D1").Value = "merk"
$ws.Range("E1").Value = "type"
$ws.Range("F1").Value = "cylinder"
$ws.Range("G1").Value = "warna"
$ws.Range("H1").Value = "no_rangka"
$ws.Range("I1").Value = "no_mesin"
$ws.Range("J1").Value = "thn_pembuatan"
$ws.Range("K1").Value = "thn_pembelian"
$ws.Range("L1").Value = "no_polisi"
$ws.Range("M1").Value = "tgl_bpkb"
$ws.Range("N1").Value = "no_bpkb"
$ws.Range("O1").Value = "harga"
$ws.Range("P1").Value = "keterangan"

# --- Row 2: sample data -------------------------------------------------
$ws.Range("A2").Value = "Tersedia"
$ws.Range("B2").Value = "Mobil"
$ws.Range("C2").Value = "17/12/2023"
$ws.Range("D2").Value = "Toyota"
$ws.Range("E2").Value = "Xenia"
$ws.Range("F2").Value = 2100
$ws.Range("G2").Value = "Hitam"
$ws.Range("H2").Value = "XENIA1023"
$ws.Range("I2").Value = "XNIA002"
$ws.Range("J2").Value = 2021
$ws.Range("K2").Value = 2021
$ws.Range("L2").Value = "B 5543 AC"

# M2 used to be a real date (numFmtId 14). It is now plain text, so the
# number-format needs clearing before the string is written, otherwise it
# would keep rendering as a date.
$ws.Range("M2").ClearFormats()
$ws.Range("M2").Value = "22/11/2021"

$ws.Range("N2").Value = "ZA22810234"
$ws.Range("O2").Value = 125000000
$ws.Range("P2").Value = "Bantuan Bupati BMS"

# --- View state: scroll to column C, leave active cell on P2 -----------
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("P2").Select()
